$d = $word.ActiveDocument

$d.Content.Find.Execute("64×43=2752", $true, $false, $false, $false, $false, $true, 1, $false, "81×90=7290", 2) | Out-Null
$d.Content.Find.Execute("13×24=312", $true, $false, $false, $false, $false, $true, 1, $false, "90×12=1080", 2) | Out-Null
$d.Content.Find.Execute("73×58=4234", $true, $false, $false, $false, $false, $true, 1, $false, "87×93=8091", 2) | Out-Null
$d.Content.Find.Execute("71×34=2414", $true, $false, $false, $false, $false, $true, 1, $false, "98×83=8134", 2) | Out-Null
$d.Content.Find.Execute("75×37=2775", $true, $false, $false, $false, $false, $true, 1, $false, "21×92=1932", 2) | Out-Null
$d.Content.Find.Execute("86×82=7052", $true, $false, $false, $false, $false, $true, 1, $false, "15×46=690", 2) | Out-Null
$d.Content.Find.Execute("63×11=693", $true, $false, $false, $false, $false, $true, 1, $false, "83×63=5229", 2) | Out-Null
$d.Content.Find.Execute("66×84=5544", $true, $false, $false, $false, $false, $true, 1, $false, "50×32=1600", 2) | Out-Null
$d.Content.Find.Execute("37×63=2331", $true, $false, $false, $false, $false, $true, 1, $false, "65×61=3965", 2) | Out-Null
$d.Content.Find.Execute("48×24=1152", $true, $false, $false, $false, $false, $true, 1, $false, "95×43=4085", 2) | Out-Null
$d.Content.Find.Execute("36×21=756", $true, $false, $false, $false, $false, $true, 1, $false, "77×35=2695", 2) | Out-Null
$d.Content.Find.Execute("58×48=2784", $true, $false, $false, $false, $false, $true, 1, $false, "46×11=506", 2) | Out-Null
$d.Content.Find.Execute("13×41=533", $true, $false, $false, $false, $false, $true, 1, $false, "36×51=1836", 2) | Out-Null
$d.Content.Find.Execute("23×91=2093", $true, $false, $false, $false, $false, $true, 1, $false, "37×66=2442", 2) | Out-Null
$d.Content.Find.Execute("40×39=1560", $true, $false, $false, $false, $false, $true, 1, $false, "94×29=2726", 2) | Out-Null
$d.Content.Find.Execute("95×85=8075", $true, $false, $false, $false, $false, $true, 1, $false, "67×15=1005", 2) | Out-Null
$d.Content.Find.Execute("30×76=2280", $true, $false, $false, $false, $false, $true, 1, $false, "89×62=5518", 2) | Out-Null
$d.Content.Find.Execute("23×98=2254", $true, $false, $false, $false, $false, $true, 1, $false, "63×54=3402", 2) | Out-Null
$d.Content.Find.Execute("73×96=7008", $true, $false, $false, $false, $false, $true, 1, $false, "59×20=1180", 2) | Out-Null
$d.Content.Find.Execute("66×89=5874", $true, $false, $false, $false, $false, $true, 1, $false, "87×40=3480", 2) | Out-Null
$d.Content.Find.Execute("69×96=6624", $true, $false, $false, $false, $false, $true, 1, $false, "32×87=2784", 2) | Out-Null
$d.Content.Find.Execute("84×49=4116", $true, $false, $false, $false, $false, $true, 1, $false, "60×34=2040", 2) | Out-Null
$d.Content.Find.Execute("97×49=4753", $true, $false, $false, $false, $false, $true, 1, $false, "69×94=6486", 2) | Out-Null
$d.Content.Find.Execute("17×12=204", $true, $false, $false, $false, $false, $true, 1, $false, "78×19=1482", 2) | Out-Null
$d.Content.Find.Execute("77×24=1848", $true, $false, $false, $false, $false, $true, 1, $false, "84×83=6972", 2) | Out-Null
